$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 43 (shifts existing rows 43:76 down to 44:77)
$ws.Rows(43).Insert()

# Seed the new row 43 with a copy of the row now at 44 (the former row 43),
# so all fields besides Fecha/Volumen/Origen start identical to the record
# it was derived from.
$ws.Range("A44:T44").Copy($ws.Range("A43:T43"))

# Apply the new record's specific values
$ws.Range("D43").Value = 44589
$ws.Range("M43").Value = 250
$ws.Range("R43").Value = "Provincia de Linares"
